$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$shp = $s.Shapes.Item(2)

# 1) Resize / reposition the content placeholder.
#    Shape.Left/Top/Width/Height are expressed in points, while the target
#    geometry is given in EMU (1 pt = 12700 EMU). A tiny epsilon is added
#    before converting so the point -> EMU truncation inside the host lands
#    on the exact EMU value instead of one below it.
function EmuToPt([double]$emu) {
    return ($emu / 12700.0) + (0.5 / 12700.0)
}
$shp.Left = EmuToPt 1451579
$shp.Top = EmuToPt 1853754
$shp.Width = EmuToPt 9603275
$shp.Height = EmuToPt 4199727

$tr = $shp.TextFrame.TextRange

# 2) Merge the three paragraphs
#      "Our goal is a model that is both accurate and generalizable. "
#      "We want this accuracy, but limit it to try to achieve the sweet spot."
#      "How?"
#    into a single paragraph.
# Setting .Text directly against the existing text would keep the common
# prefix as its own run (host preserves formatting of unchanged text), so
# the paragraph is first replaced with unrelated text (no shared prefix)
# and then with the final wording, which yields a single clean run.
$goalPara = $tr.Paragraphs(3, 1)
$goalPara.Text = "zzz"
$goalPara = $tr.Paragraphs(3, 1)
$goalPara.Text = "Our goal is a model that is both accurate and generalizable - How?"

$sweetSpotPara = $tr.Paragraphs(4, 1)
$sweetSpotPara.Delete()

$howPara = $tr.Paragraphs(4, 1)
$howPara.Delete()

# 3) Append the new paragraphs at the end of the text body.
$tail = $tr.InsertAfter("`rThe training and testing scores can reveal some info. `rTraining scores that get very accurate are an indication of overfitting. ")

$count = $tr.Paragraphs().Count
$infoPara = $tr.Paragraphs($count - 1, 1)
$infoPara.IndentLevel = 1

$run1 = $tr.InsertAfter("`rTesting scores that drop substantially (relatively) from training ")
$run2 = $run1.InsertAfter("are another. ")
